$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: sCs -> FAPs edge, values updated (target cluster stays "FAPs") ----
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1192916666666667
$ws.Range("H2").Value = 0.357875
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.558821666666667
$ws.Range("N2").Value = 7.676465
$ws.Range("O2").Value = 0.2156728774407755
$ws.Range("P2").Value = 0.2156728774407755
$ws.Range("Q2").Value = 0.3052461013194445
$ws.Range("R2").Value = 2.747214911875
$ws.Range("S2").Value = 0.2156728774407755
$ws.Range("T2").Value = 0.2156728774407755

# ---- Row 3 (new row): sCs -> ECs edge ----
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Wnt10a"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1192916666666667
$ws.Range("H3").Value = 0.357875
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("N3").Value = 19.361954
$ws.Range("O3").Value = 0.5439806384912759
$ws.Range("P3").Value = 0.5439806384912759
$ws.Range("Q3").Value = 0.7699065875277779
$ws.Range("R3").Value = 6.92915928775
$ws.Range("S3").Value = 0.5439806384912759
$ws.Range("T3").Value = 0.5439806384912759

# ---- Row 4 (new row): sCs -> sCs edge (previously row 3, with new values) ----
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Wnt10a"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1192916666666667
$ws.Range("H4").Value = 0.357875
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.851558333333334
$ws.Range("N4").Value = 8.554675000000001
$ws.Range("O4").Value = 0.2403464840679487
$ws.Range("P4").Value = 0.2403464840679487
$ws.Range("Q4").Value = 0.3401671461805556
$ws.Range("R4").Value = 3.061504315625001
$ws.Range("S4").Value = 0.2403464840679487
$ws.Range("T4").Value = 0.2403464840679487
